$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 ("bilal") entirely; row 5 ("Ahmed") shifts up to become row 4.
$ws.Rows.Item(4).Delete()

# Ensure the numeric columns (age, salary) on the new row 4 are stored as real numbers.
$ws.Cells.Item(4, 7).Value = 20
$ws.Cells.Item(4, 9).Value = 250000
